# Update Betfair back/lay odds figures for 2025-12-28 games (Sheet1).
# Each line below updates a single odds cell to its refreshed value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 1.92
$ws.Range("F3").Value = 28
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 1.13
$ws.Range("I3").Value = 1.17
$ws.Range("J3").Value = 8.4
$ws.Range("N3").Value = 5.1
$ws.Range("P3").Value = 2.42
$ws.Range("Q3").Value = 1.54
$ws.Range("R3").Value = 1.57
$ws.Range("S3").Value = 2.34
$ws.Range("U3").Value = 1.53
$ws.Range("V3").Value = 6.4
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 2.62
$ws.Range("Q4").Value = 2.1
$ws.Range("S4").Value = 3.5
$ws.Range("V4").Value = 1.4
$ws.Range("W4").Value = 1.51
$ws.Range("AC5").Value = 12
$ws.Range("AI5").Value = 36
$ws.Range("AN5").Value = 130
$ws.Range("I5").Value = 1.52
$ws.Range("L5").Value = 1.27
$ws.Range("N5").Value = 4.6
$ws.Range("P5").Value = 2.24
$ws.Range("Q5").Value = 1.66
$ws.Range("R5").Value = 1.48
$ws.Range("S5").Value = 2.66
$ws.Range("U5").Value = 1.95
$ws.Range("V5").Value = 2.92
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 12
$ws.Range("F7").Value = 1.85
$ws.Range("G7").Value = 1.9
$ws.Range("J7").Value = 3.7
$ws.Range("K7").Value = 3.9
$ws.Range("Q7").Value = 1.95
$ws.Range("T7").Value = 1.83
$ws.Range("AB8").Value = 7.4
$ws.Range("AE8").Value = 180
$ws.Range("AI8").Value = 160
$ws.Range("AM8").Value = 210
$ws.Range("J8").Value = 5.1
$ws.Range("K8").Value = 5.2
$ws.Range("N8").Value = 4.1
$ws.Range("P8").Value = 2.02
$ws.Range("Q8").Value = 1.94
$ws.Range("R8").Value = 1.39
$ws.Range("S8").Value = 3.4
$ws.Range("T8").Value = 2.24
$ws.Range("W8").Value = 3.35
$ws.Range("I9").Value = 3.6
$ws.Range("J9").Value = 3.4
$ws.Range("O9").Value = 1.37
$ws.Range("W9").Value = 1.69
$ws.Range("F10").Value = 3.8
$ws.Range("G10").Value = 4.2
$ws.Range("H10").Value = 2.08
$ws.Range("J10").Value = 3.35
$ws.Range("O10").Value = 1.37
$ws.Range("Q10").Value = 1.96
$ws.Range("V10").Value = 1.83
$ws.Range("W10").Value = 1.32
$ws.Range("F11").Value = 6.8
$ws.Range("G11").Value = 7
$ws.Range("K11").Value = 4.3
$ws.Range("P11").Value = 1.81
$ws.Range("T11").Value = 2.2
$ws.Range("U11").Value = 1.77
$ws.Range("V11").Value = 2.56
$ws.Range("AA12").Value = 48
$ws.Range("AB12").Value = 8.6
$ws.Range("AD12").Value = 14
$ws.Range("AE12").Value = 38
$ws.Range("AF12").Value = 16
$ws.Range("AG12").Value = 13.5
$ws.Range("AH12").Value = 22
$ws.Range("AI12").Value = 65
$ws.Range("AJ12").Value = 46
$ws.Range("AK12").Value = 38
$ws.Range("AL12").Value = 65
$ws.Range("AM12").Value = 150
$ws.Range("AN12").Value = 42
$ws.Range("AO12").Value = 44
$ws.Range("F12").Value = 2.84
$ws.Range("G12").Value = 2.86
$ws.Range("I12").Value = 2.9
$ws.Range("L12").Value = 1.54
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 2.88
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 1.63
$ws.Range("Q12").Value = 2.54
$ws.Range("R12").Value = 1.23
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = 2.08
$ws.Range("U12").Value = 1.88
$ws.Range("W12").Value = 1.53
$ws.Range("X12").Value = 9.6
$ws.Range("Y12").Value = 8.800000000000001
$ws.Range("Z12").Value = 16.5
$ws.Range("F13").Value = 4.4
$ws.Range("H13").Value = 2.02
$ws.Range("J13").Value = 3.45
$ws.Range("K13").Value = 3.5
$ws.Range("R13").Value = 1.25
$ws.Range("W13").Value = 1.27
$ws.Range("AB14").Value = 10.5
$ws.Range("AF14").Value = 25
$ws.Range("AJ14").Value = 80
$ws.Range("AO14").Value = 38
$ws.Range("F14").Value = 3.45
$ws.Range("G14").Value = 3.65
$ws.Range("H14").Value = 2.42
$ws.Range("J14").Value = 3.2
$ws.Range("K14").Value = 3.25
$ws.Range("N14").Value = 2.7
$ws.Range("P14").Value = 1.59
$ws.Range("Q14").Value = 2.58
$ws.Range("S14").Value = 5.3
$ws.Range("V14").Value = 1.66
$ws.Range("Z14").Value = 17.5
$ws.Range("AK15").Value = 25
$ws.Range("G15").Value = 2.16
$ws.Range("H15").Value = 3.65
$ws.Range("I15").Value = 4.7
$ws.Range("N15").Value = 4.1
$ws.Range("Q15").Value = 1.77
$ws.Range("R15").Value = 1.42
$ws.Range("W15").Value = 1.86
$ws.Range("X15").Value = 22
$ws.Range("Y15").Value = 20
$ws.Range("AB16").Value = 9.6
$ws.Range("AJ16").Value = 29
$ws.Range("F16").Value = 2.3
$ws.Range("I16").Value = 3.65
$ws.Range("P16").Value = 1.91
$ws.Range("R16").Value = 1.36
$ws.Range("W16").Value = 1.76
$ws.Range("H17").Value = 5.4
$ws.Range("Q17").Value = 2.06
$ws.Range("AA18").Value = 28
$ws.Range("AC18").Value = 7.8
$ws.Range("AE18").Value = 25
$ws.Range("AF18").Value = 26
$ws.Range("AG18").Value = 16
$ws.Range("AH18").Value = 19.5
$ws.Range("AI18").Value = 42
$ws.Range("AL18").Value = 65
$ws.Range("AM18").Value = 130
$ws.Range("AN18").Value = 55
$ws.Range("AO18").Value = 20
$ws.Range("F18").Value = 3.7
$ws.Range("G18").Value = 3.85
$ws.Range("H18").Value = 2.16
$ws.Range("I18").Value = 2.22
$ws.Range("J18").Value = 3.5
$ws.Range("V18").Value = 1.81
$ws.Range("W18").Value = 1.35
$ws.Range("X18").Value = 12
$ws.Range("Z18").Value = 13
$ws.Range("G19").Value = 4.1
$ws.Range("H19").Value = 2.1
$ws.Range("I19").Value = 2.12
$ws.Range("J19").Value = 3.65
$ws.Range("K19").Value = 3.7
$ws.Range("Q19").Value = 1.83
$ws.Range("V19").Value = 1.89
$ws.Range("H20").Value = 24
$ws.Range("J20").Value = 9.6
$ws.Range("Q20").Value = 1.42
$ws.Range("S20").Value = 2.04
$ws.Range("U20").Value = 1.67
